$d = $word.ActiveDocument

# 1. Update Ativação date: 2020 -> 2025
$d.Content.Find.Execute(
    "Ativação: 01/01/2020", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2)

# 2. Append sentence to the end of the Portuguese "Programa" (detailed) paragraph.
$d.Content.Find.Execute(
    "Movimento da Água e de solutos no Solo. Aula prática de campo: Descrição de perfil no campo. Aula prática de laboratório: Caracterização e métodos de determinação de atributos físicos e hídricos do solo.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Movimento da Água e de solutos no Solo. Aula prática de campo: Descrição de perfil no campo. Aula prática de laboratório: Caracterização e métodos de determinação de atributos físicos e hídricos do solo. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina.",
    2)

# 3. Append sentence to the end of the English (italic) "Programa" (detailed) paragraph.
$d.Content.Find.Execute(
    "Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil. The discipline may have didactic trips to complement the content of the discipline.",
    2)

# 4. Replace the "Método:" run text with the new evaluation criteria text.
$d.Content.Find.Execute(
    "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "O aluno poderá optar por um dos dois critérios de avaliação para a NF (nota final).  Critério 1: NF = média obtida em todas atividades desenvolvidas, trabalhos e relatórios ao longo do semestre. Critério 2 (alternativo): NF = (P1+P2)/2, sendo P1 e P2 avaliações escritas individuais.",
    2)

# 5. Replace the "Critério:" run text (old text -> Exame Final text with 5,0).
$d.Content.Find.Execute(
    "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.",
    2)

# 6. Replace the "Norma de recuperação:" run text (6,5 -> 5,0).
$d.Content.Find.Execute(
    "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.",
    2)
